$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# PID, Preferences, Robot now check out (were "X", now "check")
$ws.Range("C5").Value = "check"
$ws.Range("C6").Value = "check"
$ws.Range("C7").Value = "check"

# RobotMap: re-enter the same value so the formula below breaks off from
# the shared-formula group (matches the standalone <f> in the diff)
$ws.Range("C8").Value = "check"
$ws.Range("D8").Formula = '=IF(EXACT(LOWER(C8), "check"), "ü", "û")'

# ManipulatorManualControl regressed - no longer checks out
$ws.Range("C22").Value = "X"

# SyncPreferences now checks out
$ws.Range("C28").Value = "check"

# TestPID gets a note on how to exercise it
$ws.Range("E29").Value = "command to test the PID loops"

# ManipulatorManualControl note updated (added after E29's note so the
# shared-string table gets the two new entries in the same order as the
# recorded edit)
$ws.Range("E22").Value = "elevator works, arm untested"

# move the selection like the robot driver did
$ws.Range("M16").Select()
